$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$tl = $s.TimeLine
$seq = $tl.MainSequence
$eff1 = $seq.Item(6)
$eff1.MoveTo(1.0)
for ($i = 1; $i -le $seq.Count; $i++) {
    $eff2 = $seq.Item($i)
    Write-Host $i $eff2.Shape.Id $eff2.DisplayName
}
